# Adds the 37 new "Mailbox Emails" rows (A2718:D2754) that arrived in the
# SharePoint-synced mailbox export, growing Table1 to match, exactly as
# captured by the source XML diff (dimension/table ref A1:D2717 -> A1:D2754).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# New rows to append, in sheet order: DateTimeReceived, Subject, Address, HasAttachments.
# Each row is pushed with the unary "," so the jagged array of rows is not flattened.
$newData = @()
$newData += ,@('2025-07-25T18:17:37+00:00', 'EXTERNAL:- BBC Radio 2 - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:44+00:00', 'EXTERNAL:- BBC Radio Shetland - Wk30 - 2025-07-29 - Tuesday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:42+00:00', 'EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-29 - Tuesday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:40+00:00', 'EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-28 - Monday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:39+00:00', 'EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-30 - Wednesday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:38+00:00', 'EXTERNAL:- BBC World Service UK Schedule - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:37+00:00', 'EXTERNAL:- BBC Radio 6 Music - Wk30 - 2025-07-28 - Monday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:36+00:00', 'EXTERNAL:- BBC World Service UK Schedule - Wk29 - 2025-07-25 - Friday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:32+00:00', 'EXTERNAL:- BBC Radio 1 - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:30+00:00', 'EXTERNAL:- BBC Asian Network - Wk31 - 2025-08-03 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:28+00:00', 'EXTERNAL:- BBC Radio nan Gàidheal - Wk30 - 2025-07-29 - Tuesday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:25+00:00', 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-28 - Monday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:23+00:00', 'EXTERNAL:- BBC Asian Network - Wk31 - 2025-08-02 - Saturday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:20+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-08-01 - Friday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:19+00:00', 'EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-31 - Thursday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:15+00:00', 'EXTERNAL:- BBC Radio nan Gàidheal - Wk30 - 2025-07-31 - Thursday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:13+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:12+00:00', 'EXTERNAL:- BBC Radio 6 Music - Wk30 - 2025-07-30 - Wednesday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:11+00:00', 'EXTERNAL:- BBC Radio 2 - Wk31 - 2025-08-03 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:09+00:00', 'EXTERNAL:- BBC Radio 4 FM - Wk30 - 2025-07-31 - Thursday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:09+00:00', 'EXTERNAL:- BBC World Service UK Schedule - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:09+00:00', 'EXTERNAL:- BBC Radio Shetland - Wk30 - 2025-07-28 - Monday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:08+00:00', 'EXTERNAL:- BBC Radio Cymru 2 - Wk30 - 2025-07-28 - Monday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:04+00:00', 'EXTERNAL:- BBC Radio 6 Music - Wk30 - 2025-07-29 - Tuesday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:04+00:00', 'EXTERNAL:- BBC Four HD - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:04+00:00', 'EXTERNAL:- BBC Four HD - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:03+00:00', 'EXTERNAL:- BBC Radio Cymru - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:00+00:00', 'EXTERNAL:- BBC Radio 6 Music - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:16:00+00:00', 'EXTERNAL:- BBC Radio Shetland - Wk30 - 2025-07-31 - Thursday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:58+00:00', 'EXTERNAL:- BBC Asian Network - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:58+00:00', 'EXTERNAL:- BBC Radio 1Xtra - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:58+00:00', 'EXTERNAL:- BBC Radio 6 Music - Wk30 - 2025-07-27 - Sunday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:56+00:00', 'EXTERNAL:- BBC Radio nan Gàidheal - Wk30 - 2025-08-01 - Friday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:56+00:00', 'EXTERNAL:- BBC Radio 1Xtra - Wk31 - 2025-08-04 - Monday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:55+00:00', 'EXTERNAL:- BBC Radio 1Xtra - Wk30 - 2025-08-01 - Friday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:53+00:00', 'EXTERNAL:- BBC Radio 2 - Wk30 - 2025-07-26 - Saturday', 'pressportal@bbc.co.uk', $true)
$newData += ,@('2025-07-25T18:15:52+00:00', 'EXTERNAL:- BBC Radio 4 Extra - Wk31 - 2025-08-02 - Saturday', 'pressportal@bbc.co.uk', $true)

# Grow the table by one ListRow per new record so Table1 + autoFilter + the
# sheet dimension all expand together (mirrors Excel UI "Insert Table Row").
$startRow = $lo.Range.Row + $lo.ListRows.Count + 1
foreach ($rowData in $newData) {
    $lo.ListRows.Add() | Out-Null
}

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
